{"js": "// Replace the date line and all 25 two-digit multiplication problems in\n// the document body with their new values (per the commit's diff).\nconst replacements = [\n  [\"2024-06-08 Saturday\", \"2024-06-09 Sunday\"],\n  [\"63\u00d769=\", \"16\u00d740=\"],\n  [\"32\u00d745=\", \"20\u00d781=\"],\n  [\"77\u00d769=\", \"60\u00d729=\"],\n  [\"34\u00d784=\", \"69\u00d761=\"],\n  [\"34\u00d795=\", \"94\u00d765=\"],\n  [\"74\u00d730=\", \"52\u00d799=\"],\n  [\"42\u00d729=\", \"92\u00d756=\"],\n  [\"75\u00d781=\", \"38\u00d772=\"],\n  [\"93\u00d766=\", \"91\u00d728=\"],\n  [\"42\u00d749=\", \"70\u00d734=\"],\n  [\"72\u00d758=\", \"58\u00d716=\"],\n  [\"89\u00d768=\", \"28\u00d736=\"],\n  [\"39\u00d786=\", \"54\u00d762=\"],\n  [\"75\u00d723=\", \"46\u00d726=\"],\n  [\"81\u00d791=\", \"39\u00d793=\"],\n  [\"62\u00d785=\", \"61\u00d737=\"],\n  [\"43\u00d745=\", \"30\u00d751=\"],\n  [\"54\u00d753=\", \"71\u00d766=\"],\n  [\"83\u00d758=\", \"89\u00d719=\"],\n  [\"63\u00d747=\", \"29\u00d731=\"],\n  [\"82\u00d733=\", \"72\u00d727=\"],\n  [\"53\u00d736=\", \"45\u00d746=\"],\n  [\"49\u00d751=\", \"98\u00d798=\"],\n  [\"80\u00d762=\", \"62\u00d755=\"],\n  [\"83\u00d734=\", \"80\u00d743=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all 25 two-digit multiplication problems to\n# their new values (per the commit's diff), using Find/Replace over the\n# whole document content.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-06-08 Saturday\", \"2024-06-09 Sunday\"),\n  @(\"63\u00d769=\", \"16\u00d740=\"),\n  @(\"32\u00d745=\", \"20\u00d781=\"),\n  @(\"77\u00d769=\", \"60\u00d729=\"),\n  @(\"34\u00d784=\", \"69\u00d761=\"),\n  @(\"34\u00d795=\", \"94\u00d765=\"),\n  @(\"74\u00d730=\", \"52\u00d799=\"),\n  @(\"42\u00d729=\", \"92\u00d756=\"),\n  @(\"75\u00d781=\", \"38\u00d772=\"),\n  @(\"93\u00d766=\", \"91\u00d728=\"),\n  @(\"42\u00d749=\", \"70\u00d734=\"),\n  @(\"72\u00d758=\", \"58\u00d716=\"),\n  @(\"89\u00d768=\", \"28\u00d736=\"),\n  @(\"39\u00d786=\", \"54\u00d762=\"),\n  @(\"75\u00d723=\", \"46\u00d726=\"),\n  @(\"81\u00d791=\", \"39\u00d793=\"),\n  @(\"62\u00d785=\", \"61\u00d737=\"),\n  @(\"43\u00d745=\", \"30\u00d751=\"),\n  @(\"54\u00d753=\", \"71\u00d766=\"),\n  @(\"83\u00d758=\", \"89\u00d719=\"),\n  @(\"63\u00d747=\", \"29\u00d731=\"),\n  @(\"82\u00d733=\", \"72\u00d727=\"),\n  @(\"53\u00d736=\", \"45\u00d746=\"),\n  @(\"49\u00d751=\", \"98\u00d798=\"),\n  @(\"80\u00d762=\", \"62\u00d755=\"),\n  @(\"83\u00d734=\", \"80\u00d743=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  [void]$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2)\n}\n"}
